$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.25146210193634
$ws.Range("B1").Value = 2.489693164825439
$ws.Range("C1").Value = 4.887235164642334
$ws.Range("D1").Value = 3.066373586654663
$ws.Range("E1").Value = 1.118917107582092
